$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.359437
$ws.Range("H2").Value = 55.078311
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.672264666666667
$ws.Range("N2").Value = 8.016794000000001
$ws.Range("O2").Value = 0.06772620019093417
$ws.Range("P2").Value = 0.06772620019093417
$ws.Range("Q2").Value = 49.06127479499267
$ws.Range("R2").Value = 441.5514731549341
$ws.Range("S2").Value = 0.06772620019093417
$ws.Range("T2").Value = 0.06772620019093417

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.359437
$ws.Range("H3").Value = 55.078311
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.41886
$ws.Range("N3").Value = 82.25658
$ws.Range("O3").Value = 0.6949069171668364
$ws.Range("P3").Value = 0.6949069171668364
$ws.Range("Q3").Value = 503.39483278182
$ws.Range("R3").Value = 4530.55349503638
$ws.Range("S3").Value = 0.6949069171668364
$ws.Range("T3").Value = 0.6949069171668364

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.359437
$ws.Range("H4").Value = 55.078311
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.365757
$ws.Range("N4").Value = 28.097271
$ws.Range("O4").Value = 0.2373668826422294
$ws.Range("P4").Value = 0.2373668826422294
$ws.Range("Q4").Value = 171.950025598809
$ws.Range("R4").Value = 1547.550230389281
$ws.Range("S4").Value = 0.2373668826422294
$ws.Range("T4").Value = 0.2373668826422294
